$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 90.56466633333334
$ws.Range("H2").Value = 271.693999
$ws.Range("I2").Value = 0.2234788625831797
$ws.Range("J2").Value = 0.2234788625831796
$ws.Range("Q2").Value = 1.347783364372667
$ws.Range("R2").Value = 12.130050279354
$ws.Range("S2").Value = 0.2234788625831797
$ws.Range("T2").Value = 0.2234788625831796

# Row 3
$ws.Range("I3").Value = 0.601197186834308
$ws.Range("J3").Value = 0.6011971868343079
$ws.Range("S3").Value = 0.601197186834308
$ws.Range("T3").Value = 0.6011971868343079

# Row 4
$ws.Range("I4").Value = 0.1753239505825123
$ws.Range("J4").Value = 0.1753239505825123
$ws.Range("S4").Value = 0.1753239505825123
$ws.Range("T4").Value = 0.1753239505825123
